# Git_Comandos_Utiles.xlsx
# Adds two new rows to "Tabla1":
#   - git config --global user.name / Estado / Devuelve el Nombre de Usuario de GIT
#   - git status                    / Estado / Estado de los archivos
# then re-sorts the table by the "Clasificacion" column (matches the
# sortCondition the author applied), resizes the table range to A1:C12
# (author resized/left the table one row taller than the data, as in the
# source workbook) and moves the active selection to A7.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- add the new command rows -------------------------------------------
# Order matters: added in the same order the sort leaves them in when two
# rows share the same "Clasificacion" key (stable sort keeps insertion
# order), so "git config --global user.name" must be added before
# "git status".
$rowConfig = $tbl.ListRows.Add()
$rowConfig.Range.Item(1, 1).Value = "git config --global user.name"
$rowConfig.Range.Item(1, 2).Value = "Estado"
$rowConfig.Range.Item(1, 3).Value = "Devuelve el Nombre de Usuario de GIT"

$rowStatus = $tbl.ListRows.Add()
$rowStatus.Range.Item(1, 1).Value = "git status"
$rowStatus.Range.Item(1, 2).Value = "Estado"
$rowStatus.Range.Item(1, 3).Value = "Estado de los archivos"

# --- sort the table by "Clasificacion" (column B) -----------------------
$sortObj = $tbl.Sort
$sortObj.SortFields.Clear()
$colClasificacion = $tbl.ListColumns.Item("Clasificacion").Range
$sortObj.SortFields.Add($colClasificacion)
$sortObj.Header = 1
$sortObj.Apply()

# --- resize the table range to match the author's saved extent ----------
$tbl.Resize($ws.Range("A1:C12"))

# --- restore the active-cell selection -----------------------------------
$ws.Range("A7").Select()
